# Rename "PRISMA carlona" -> "PRISMA"
$wb = $excel.ActiveWorkbook
$wsPrisma = $wb.Worksheets.Item("PRISMA carlona")
$wsPrisma.Name = "PRISMA"

# Add new sheet "RETICOLO" right after "PRISMA"
$wsReticolo = $wb.Worksheets.Add($null, $wsPrisma)
$wsReticolo.Name = "RETICOLO"

# Match the page/print setup used by the rest of the workbook
$wsReticolo.PageSetup.PaperSize = $wsPrisma.PageSetup.PaperSize
$wsReticolo.PageSetup.Zoom = $wsPrisma.PageSetup.Zoom
$wsReticolo.PageSetup.FitToPagesWide = $wsPrisma.PageSetup.FitToPagesWide
$wsReticolo.PageSetup.FitToPagesTall = $wsPrisma.PageSetup.FitToPagesTall
$wsReticolo.PageSetup.Orientation = $wsPrisma.PageSetup.Orientation
$wsReticolo.PageSetup.LeftMargin = $wsPrisma.PageSetup.LeftMargin
$wsReticolo.PageSetup.RightMargin = $wsPrisma.PageSetup.RightMargin
$wsReticolo.PageSetup.TopMargin = $wsPrisma.PageSetup.TopMargin
$wsReticolo.PageSetup.BottomMargin = $wsPrisma.PageSetup.BottomMargin
$wsReticolo.PageSetup.HeaderMargin = $wsPrisma.PageSetup.HeaderMargin
$wsReticolo.PageSetup.FooterMargin = $wsPrisma.PageSetup.FooterMargin
$wsReticolo.PageSetup.PrintHeadings = $wsPrisma.PageSetup.PrintHeadings
$wsReticolo.PageSetup.PrintGridlines = $wsPrisma.PageSetup.PrintGridlines
$wsReticolo.PageSetup.BlackAndWhite = $wsPrisma.PageSetup.BlackAndWhite
$wsReticolo.PageSetup.Draft = $wsPrisma.PageSetup.Draft
$wsReticolo.PageSetup.Order = $wsPrisma.PageSetup.Order
$wsReticolo.DisplayGridlines = $wsPrisma.DisplayGridlines

# Match the sheet's zoom level (100%) like the rest of the workbook
$wsReticolo.Activate()
$wsReticolo.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100

# Make "ERRORI" the active sheet again (activeTab = 0)
$wb.Worksheets.Item("ERRORI").Activate()
